$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 251044.17
$ws.Range("I15").Value = 251044.17
$ws.Range("K15").Value = 753132.51
$ws.Range("M15").Value = -752963.51
$ws.Range("H136").Value = 39000
$ws.Range("J136").Value = 39000
$ws.Range("L136").Value = 39000
$ws.Range("N136").Value = -49200
$ws.Range("H138").Value = 1982.61
$ws.Range("I138").Value = 662.2083
$ws.Range("J138").Value = 2399.5789
$ws.Range("K138").Value = 1986.6249
$ws.Range("L138").Value = 7198.736699999999
$ws.Range("M138").Value = 3153.3751
$ws.Range("N138").Value = -17478.7367
$ws.Range("H141").Value = 1827.069
$ws.Range("I141").Value = 1406.8518
$ws.Range("J141").Value = 7500
$ws.Range("K141").Value = 4220.555399999999
$ws.Range("L141").Value = 22500
$ws.Range("M141").Value = 959.4446000000007
$ws.Range("N141").Value = -32860

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1407.6471
$ws.Range("I2").Value = 833.75
$ws.Range("J2").Value = 1917.7778
$ws.Range("K2").Value = 833.75
$ws.Range("L2").Value = 1917.7778
$ws.Range("M2").Value = -720.75
$ws.Range("N2").Value = -2143.7778
$ws.Range("H32").Value = 27435.88
$ws.Range("I32").Value = 2424.8057
$ws.Range("J32").Value = 177502.33
$ws.Range("K32").Value = 2424.8057
$ws.Range("L32").Value = 177502.33
$ws.Range("M32").Value = -2137.8057
$ws.Range("N32").Value = -178076.33
$ws.Range("H61").Value = 3495.5173
$ws.Range("I61").Value = 1922.5834
$ws.Range("J61").Value = 4605.8237
$ws.Range("K61").Value = 1922.5834
$ws.Range("L61").Value = 4605.8237
$ws.Range("M61").Value = -1710.5834
$ws.Range("N61").Value = -5029.8237
$ws.Range("H110").Value = 615.3333
$ws.Range("I110").Value = 574.6667
$ws.Range("K110").Value = 574.6667
$ws.Range("M110").Value = 1470.3333
$ws.Range("H116").Value = 1407.6471
$ws.Range("I116").Value = 833.75
$ws.Range("J116").Value = 1917.7778
$ws.Range("K116").Value = 833.75
$ws.Range("L116").Value = 1917.7778
$ws.Range("M116").Value = 1460.25
$ws.Range("N116").Value = -6505.7778
$ws.Range("H122").Value = 2385
$ws.Range("I122").Value = 2356.25
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 7068.75
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -4618.75
$ws.Range("N122").Value = -12400
$ws.Range("H133").Value = 36919.8
$ws.Range("J133").Value = 36919.8
$ws.Range("L133").Value = 36919.8
$ws.Range("N133").Value = -41979.8
$ws.Range("H136").Value = 3495.5173
$ws.Range("I136").Value = 1922.5834
$ws.Range("J136").Value = 4605.8237
$ws.Range("K136").Value = 5767.7502
$ws.Range("L136").Value = 13817.4711
$ws.Range("M136").Value = -3217.7502
$ws.Range("N136").Value = -18917.4711
$ws.Range("H137").Value = 41500
$ws.Range("J137").Value = 41500
$ws.Range("L137").Value = 41500
$ws.Range("N137").Value = -51700

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1407.6471
$ws.Range("I3").Value = 833.75
$ws.Range("J3").Value = 1917.7778
$ws.Range("K3").Value = 833.75
$ws.Range("L3").Value = 1917.7778
$ws.Range("M3").Value = -719.75
$ws.Range("N3").Value = -2145.7778
$ws.Range("H94").Value = 899.04346
$ws.Range("I94").Value = 946.2941
$ws.Range("J94").Value = 765.1667
$ws.Range("K94").Value = 946.2941
$ws.Range("L94").Value = 765.1667
$ws.Range("M94").Value = -495.2941
$ws.Range("N94").Value = -1667.1667
$ws.Range("H107").Value = 896.9048
$ws.Range("I107").Value = 629.1429000000001
$ws.Range("J107").Value = 1432.4286
$ws.Range("K107").Value = 629.1429000000001
$ws.Range("L107").Value = 1432.4286
$ws.Range("M107").Value = 1290.8571
$ws.Range("N107").Value = -5272.4286

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1104.8541
$ws.Range("I31").Value = 875.4186
$ws.Range("K31").Value = 875.4186
$ws.Range("M31").Value = -580.4186
$ws.Range("H34").Value = 1104.8541
$ws.Range("I34").Value = 875.4186
$ws.Range("K34").Value = 875.4186
$ws.Range("M34").Value = -673.4186
$ws.Range("H99").Value = 5209630
$ws.Range("I99").Value = 6945601
$ws.Range("K99").Value = 6945601
$ws.Range("M99").Value = -6944103
$ws.Range("H126").Value = 5209630
$ws.Range("I126").Value = 6945601
$ws.Range("K126").Value = 20836803
$ws.Range("M126").Value = -20834333
$ws.Range("H132").Value = 1895.4667
$ws.Range("I132").Value = 1463.4054
$ws.Range("J132").Value = 3893.75
$ws.Range("K132").Value = 4390.216200000001
$ws.Range("L132").Value = 11681.25
$ws.Range("M132").Value = -1860.216200000001
$ws.Range("N132").Value = -16741.25
$ws.Range("H134").Value = 1961.0238
$ws.Range("I134").Value = 1017.2857
$ws.Range("J134").Value = 6679.7144
$ws.Range("K134").Value = 3051.8571
$ws.Range("L134").Value = 20039.1432
$ws.Range("M134").Value = -516.8571000000002
$ws.Range("N134").Value = -25109.1432

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 2999.5
$ws.Range("H79").Value = 2999.5
$ws.Range("H131").Value = 1362.4
$ws.Range("I131").Value = 332.7857
$ws.Range("J131").Value = 1598.705
$ws.Range("K131").Value = 998.3571000000001
$ws.Range("L131").Value = 4796.115
$ws.Range("M131").Value = 4041.6429
$ws.Range("N131").Value = -14876.115
$ws.Range("H132").Value = 1001.2083
$ws.Range("J132").Value = 1164.5
$ws.Range("L132").Value = 10480.5
$ws.Range("N132").Value = -15540.5
$ws.Range("H140").Value = 5793.077
$ws.Range("I140").Value = 10112.272
$ws.Range("J140").Value = 2625.6667
$ws.Range("K140").Value = 30336.816
$ws.Range("L140").Value = 7877.000100000001
$ws.Range("M140").Value = -25156.816
$ws.Range("N140").Value = -18237.0001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 2907.5
$ws.Range("I55").Value = 1815
$ws.Range("J55").Value = 4000
$ws.Range("K55").Value = 1815
$ws.Range("L55").Value = 4000
$ws.Range("M55").Value = -1488
$ws.Range("N55").Value = -4654
$ws.Range("H122").Value = 927166.4399999999
$ws.Range("I122").Value = 1112211.6
$ws.Range("J122").Value = 1940.5
$ws.Range("K122").Value = 3336634.8
$ws.Range("L122").Value = 5821.5
$ws.Range("M122").Value = -3334184.8
$ws.Range("N122").Value = -10721.5
$ws.Range("H137").Value = 44999
$ws.Range("J137").Value = 44999
$ws.Range("L137").Value = 44999
$ws.Range("N137").Value = -55199
$ws.Range("H138").Value = 49697
$ws.Range("J138").Value = 49697
$ws.Range("L138").Value = 49697
$ws.Range("N138").Value = -59977
$ws.Range("H139").Value = 35313
$ws.Range("J139").Value = 35313
$ws.Range("L139").Value = 35313
$ws.Range("N139").Value = -45593

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 42500
$ws.Range("I34").Value = 42500
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 42500
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -42328
$ws.Range("N34").ClearContents()
$ws.Range("H45").Value = 4999
$ws.Range("I45").Value = 4999
$ws.Range("K45").Value = 4999
$ws.Range("M45").Value = -4592
$ws.Range("H122").Value = 3239.4644
$ws.Range("I122").Value = 2254.5454
$ws.Range("J122").Value = 3876.7646
$ws.Range("K122").Value = 6763.6362
$ws.Range("L122").Value = 11630.2938
$ws.Range("M122").Value = -4313.6362
$ws.Range("N122").Value = -16530.2938
$ws.Range("H132").Value = 5493.393
$ws.Range("I132").Value = 4889.8823
$ws.Range("J132").Value = 6426.091
$ws.Range("K132").Value = 14669.6469
$ws.Range("L132").Value = 19278.273
$ws.Range("M132").Value = -12139.6469
$ws.Range("N132").Value = -24338.273

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H126").Value = 73300
$ws.Range("I126").Value = 78784.62
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 236353.86
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -233883.86
$ws.Range("N126").Value = -10940
$ws.Range("H135").Value = 116206.43
$ws.Range("I135").Value = 40000
$ws.Range("J135").Value = 128907.5
$ws.Range("K135").Value = 40000
$ws.Range("L135").Value = 128907.5
$ws.Range("M135").Value = -34930
$ws.Range("N135").Value = -139047.5
$ws.Range("H137").Value = 40000
$ws.Range("J137").Value = 40000
$ws.Range("L137").Value = 40000
$ws.Range("N137").Value = -50200
$ws.Range("H139").Value = 50000
$ws.Range("J139").Value = 50000
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280
